$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "23.423.38"
$ws.Range("E2").Value = "  +0.93%  "
$ws.Range("D3").Value = "1.638.08"
$ws.Range("E3").Value = "  +2.23%  "
$ws.Range("D4").Value = "1.004"
$ws.Range("E4").Value = "  +0.20%  "
$ws.Range("D5").Value = "1.004"
$ws.Range("E5").Value = "  +0.24%  "
$ws.Range("D6").Value = "303.62"
$ws.Range("E6").Value = "  -0.42%  "
$ws.Range("D7").Value = "0.3765"
$ws.Range("E7").Value = "  +0.10%  "
$ws.Range("D8").Value = "52.19"
$ws.Range("E8").Value = "  -1.16%  "
$ws.Range("D9").Value = "0.3637"
$ws.Range("E9").Value = "  +0.85%  "
$ws.Range("D10").Value = "1.243"
$ws.Range("E10").Value = "  -1.43%  "
$ws.Range("D11").Value = "0.08096"
$ws.Range("E11").Value = "  -0.41%  "
$ws.Range("D12").Value = "1.004"
$ws.Range("E12").Value = "  +0.25%  "
$ws.Range("D13").Value = "22.83"
$ws.Range("E13").Value = "  +0.06%  "
$ws.Range("D14").Value = "6.626"
$ws.Range("E14").Value = "  +0.62%  "
$ws.Range("D15").Value = "0.00001251"
$ws.Range("E15").Value = "  +0.65%  "
$ws.Range("D16").Value = "7.282"
$ws.Range("E16").Value = "  -0.88%  "
$ws.Range("D17").Value = "1.638.29"
$ws.Range("E17").Value = "  +2.32%  "
$ws.Range("D18").Value = "94.07"
$ws.Range("E18").Value = "  +0.09%  "
$ws.Range("D19").Value = "0.06929"
$ws.Range("E19").Value = "  +0.02%  "
$ws.Range("D20").Value = "18.12"
$ws.Range("E20").Value = "  -0.02%  "
$ws.Range("D21").Value = "6.539"
$ws.Range("E21").Value = "  +0.19%  "
$ws.Range("D22").Value = "1.003"
$ws.Range("E22").Value = "  -0.05%  "
$ws.Range("D23").Value = "23.444.71"
$ws.Range("E23").Value = "  +0.99%  "
$ws.Range("D24").Value = "12.83"
$ws.Range("E24").Value = "  -0.26%  "
$ws.Range("D25").Value = "3.226"
$ws.Range("E25").Value = "  +5.61%  "
$ws.Range("D26").Value = "2.438"
$ws.Range("E26").Value = "  +0.58%  "
$ws.Range("D27").Value = "21.15"
$ws.Range("E27").Value = "  -0.02%  "
$ws.Range("D28").Value = "150.98"
$ws.Range("E28").Value = "  +0.21%  "
$ws.Range("D29").Value = "5.313"
$ws.Range("E29").Value = "  +0.98%  "
$ws.Range("D30").Value = "135.58"
$ws.Range("E30").Value = "  +0.49%  "
$ws.Range("D31").Value = "2.312"
$ws.Range("E31").Value = "  -3.98%  "
$ws.Range("D32").Value = "1.820.21"
$ws.Range("E32").Value = "  +2.24%  "
$ws.Range("D33").Value = "6.863"
$ws.Range("E33").Value = "  +1.87%  "
$ws.Range("D34").Value = "10.97"
$ws.Range("E34").Value = "  +6.85%  "
$ws.Range("D35").Value = "0.9611"
$ws.Range("E35").Value = "  +1.02%  "
$ws.Range("D36").Value = "0.02857"
$ws.Range("E36").Value = "  +3.33%  "
$ws.Range("D37").Value = "6.253"
$ws.Range("E37").Value = "  +2.58%  "
$ws.Range("D38").Value = "0.2551"
$ws.Range("E38").Value = "  +1.55%  "
$ws.Range("D39").Value = "0.07271"
$ws.Range("E39").Value = "  -1.98%  "
$ws.Range("D40").Value = "0.08844"
$ws.Range("E40").Value = "  +1.14%  "
$ws.Range("D41").Value = "1.374"
$ws.Range("E41").Value = "  -2.00%  "
$ws.Range("D42").Value = "0.7108"
$ws.Range("E42").Value = "  +0.33%  "
$ws.Range("D43").Value = "16.36"
$ws.Range("E43").Value = "  +3.14%  "
$ws.Range("D44").Value = "12.58"
$ws.Range("E44").Value = "  +1.36%  "
$ws.Range("D45").Value = "0.6549"
$ws.Range("E45").Value = "  +0.54%  "
$ws.Range("D46").Value = "2.349"
$ws.Range("E46").Value = "  +1.18%  "
$ws.Range("D47").Value = "1.003"
$ws.Range("E47").Value = "  +0.27%  "
$ws.Range("D48").Value = "3.993"
$ws.Range("E48").Value = "  -0.44%  "
$ws.Range("D49").Value = "0.07990"
$ws.Range("E49").Value = "  +0.58%  "
$ws.Range("D50").Value = "1.216"
$ws.Range("E50").Value = "  +1.73%  "
$ws.Range("D51").Value = "127.87"
$ws.Range("E51").Value = "  -4.49%  "

$ws.Range("D2:D51").ClearFormats()
